{"js": "// Apply the English -> Portuguese translation edits described by the diff.\n// Strategy: load all body paragraphs once, then for each paragraph that\n// needs a change, run a paragraph-scoped search() (so duplicate phrases in\n// other paragraphs are left untouched) and insertText(..., \"Replace\") on\n// each hit. This preserves existing run formatting (bold, highlight,\n// hyperlink, color, etc.) because insertText(\"Replace\") only rewrites the\n// matched text, not the surrounding run properties.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Replace `needle` with `replacement` inside paragraph `idx`. Expects exactly\n// `expectedCount` matches (defaults to 1) - helps catch drift early.\nasync function replaceInParagraph(idx, needle, replacement, expectedCount) {\n  const para = paragraphs.items[idx];\n  const results = para.search(needle, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  const count = results.items.length;\n  if (expectedCount !== undefined && count !== expectedCount) {\n    throw new Error(\n      `Paragraph ${idx}: expected ${expectedCount} match(es) for \"${needle}\", found ${count}`\n    );\n  }\n  for (let i = 0; i < count; i++) {\n    results.items[i].insertText(replacement, \"Replace\");\n  }\n  await context.sync();\n}\n\n// Paragraph 0: \"English\" link label + sibling language list.\nawait replaceInParagraph(0, \"English\", \"Ingl\u00eas\", 1);\nawait replaceInParagraph(\n  0,\n  \" / Portuguese / French / Thai / Vietnamese / Spanish\",\n  \" / Portugu\u00eas / Franc\u00eas / Tailand\u00eas / Vietnamita / Espanhol\",\n  1\n);\n\n// Paragraph 2: secondary \"English\" label (below the bookmark).\nawait replaceInParagraph(2, \"English\", \"Ingl\u00eas\", 1);\n\n// Paragraph 4: \"Brief\" bold label.\nawait replaceInParagraph(4, \"Brief\", \"Resumo\", 1);\n\n// Paragraph 5: brief description text.\nawait replaceInParagraph(\n  5,\n  \"An email sent to partners in the target country who RSVPed yes but haven\\u2019t sent their documents to us. It will be sent via customer.io\",\n  \"E-mail enviado a parceiros no pa\\u00EDs-alvo que responderam que sim, mas ainda n\\u00E3o nos enviaram os documentos. Ser\\u00E1 enviado atrav\\u00E9s do customer.io\",\n  1\n);\n\n// Paragraph 7: \"Target audience\" bold label.\nawait replaceInParagraph(7, \"Target audience\", \"P\\u00FAblico-alvo\", 1);\n\n// Paragraph 8: target audience description.\nawait replaceInParagraph(\n  8,\n  \"Invited partners who haven\\u2019t submitted their documents\",\n  \"Parceiros convidados que ainda n\\u00E3o enviaram os seus documentos\",\n  1\n);\n\n// Paragraph 10: first \"Subject line\" block.\nawait replaceInParagraph(10, \"Subject line\", \"Linha de assunto\", 1);\nawait replaceInParagraph(10, \"[EVENT NAME]\", \"[NOME DO EVENTO]\", 1);\nawait replaceInParagraph(\n  10,\n  \" \\u2014 have you submitted your docs?  \",\n  \" - j\\u00E1 enviou os seus documentos?  \",\n  1\n);\n\n// Paragraph 12: \"Don't forget to send your documents\" headline (block 1).\nawait replaceInParagraph(\n  12,\n  \"Don\\u2019t forget to send your documents\",\n  \"N\\u00E3o se esque\\u00E7a de enviar os seus documentos\",\n  1\n);\n\n// Paragraph 14: \"Hi [PARTNER NAME], \" greeting (block 1).\nawait replaceInParagraph(14, \"Hi \", \"Ol\\u00E1 \", 1);\nawait replaceInParagraph(14, \"[PARTNER NAME]\", \"[NOME DO PARCEIRO]\", 1);\n\n// Paragraph 16: \"We're excited...\" intro (block 1).\nawait replaceInParagraph(\n  16,\n  \"We\\u2019re excited to see you at the upcoming \",\n  \"Estamos ansiosos por v\\u00EA-lo no pr\\u00F3ximo evento \",\n  1\n);\nawait replaceInParagraph(16, \"[EVENT NAME]\", \"[NOME DO EVENTO]\", 1);\n\n// Paragraph 17: \"To confirm your registration...\" (block 1). Note: the\n// \"DD Mmm YYYY\" placeholder that follows is NOT changed in this block.\nawait replaceInParagraph(\n  17,\n  \"To confirm your registration, we need the following documents from you by \",\n  \"Para confirmar a sua inscri\\u00E7\\u00E3o, precisamos que nos envie os seguintes documentos at\\u00E9 dia \",\n  1\n);\n\n// Paragraph 18: document list placeholder (block 1).\nawait replaceInParagraph(\n  18,\n  \"[insert list of documents required]\",\n  \"[inserir lista dos documentos necess\\u00E1rios].\",\n  1\n);\n\n// Paragraph 19: \"Please send a copy...\" (block 1).\nawait replaceInParagraph(\n  19,\n  \"Please send a copy of these documents to your country manager, \",\n  \"Por favor, envie uma c\\u00F3pia destes documentos ao seu gestor de parcerias, \",\n  1\n);\nawait replaceInParagraph(19, \", at \", \", para \", 1);\nawait replaceInParagraph(19, \" or \", \" ou \", 1);\nawait replaceInParagraph(\n  19,\n  \" (WhatsApp), so that we can make the necessary arrangements for you, including accommodation and transportation.\",\n  \" (WhatsApp), para podermos tomar as medidas necess\\u00E1rias, incluindo alojamento e transporte.\",\n  1\n);\n\n// Paragraph 20: \"If you have any questions, please contact your country manager.\" (block 1).\nawait replaceInParagraph(\n  20,\n  \"If you have any questions, please contact your country manager.\",\n  \"Se tiver alguma d\\u00FAvida, contacte o gestor do seu pa\\u00EDs.\",\n  1\n);\n\n// Paragraph 21: closing line (block 1).\nawait replaceInParagraph(\n  21,\n  \"We look forward to seeing you there!\",\n  \"Contamos com a sua presen\\u00E7a!\",\n  1\n);\n\n// Paragraph 27: second \"Subject line\" block.\nawait replaceInParagraph(27, \"Subject line\", \"Linha de assunto\", 1);\nawait replaceInParagraph(27, \"[EVENT NAME]\", \"[NOME DO EVENTO]\", 1);\nawait replaceInParagraph(\n  27,\n  \" \\u2014 have you submitted your docs?  \",\n  \" - j\\u00E1 enviou os seus documentos?  \",\n  1\n);\n\n// Paragraph 29: \"Don't forget to send your documents\" headline (block 2).\nawait replaceInParagraph(\n  29,\n  \"Don\\u2019t forget to send your documents\",\n  \"N\\u00E3o se esque\\u00E7a de enviar os seus documentos\",\n  1\n);\n\n// Paragraph 31: \"Dear [PARTNER NAME], \" greeting (block 2).\nawait replaceInParagraph(31, \"Dear \", \"Ol\\u00E1 \", 1);\nawait replaceInParagraph(31, \"[PARTNER NAME]\", \"[NOME DO PARCEIRO]\", 1);\n\n// Paragraph 33: \"We're excited...\" intro (block 2).\nawait replaceInParagraph(\n  33,\n  \"We\\u2019re excited to see you at the upcoming \",\n  \"Estamos ansiosos por v\\u00EA-lo no pr\\u00F3ximo evento \",\n  1\n);\nawait replaceInParagraph(33, \"[EVENT NAME]\", \"[NOME DO EVENTO]\", 1);\n\n// Paragraph 34: \"To ensure you have the best experience...\" (block 2) -\n// here the \"DD Mmm YYYY\" placeholder IS changed to \"DD Mmm AAAA\".\nawait replaceInParagraph(\n  34,\n  \"To ensure you have the best experience at this event, we need the following documents from you by \",\n  \"De forma a garantir a melhor experi\\u00EAncia poss\\u00EDvel neste evento, \\u00E9 necess\\u00E1rio que nos envie os seguintes documentos at\\u00E9 \",\n  1\n);\nawait replaceInParagraph(34, \"DD Mmm YYYY\", \"DD Mmm AAAA\", 1);\n\n// Paragraph 35: document list placeholder (block 2).\nawait replaceInParagraph(\n  35,\n  \"[insert list of documents required]\",\n  \"[inserir lista dos documentos necess\\u00E1rios].\",\n  1\n);\n\n// Paragraph 36: \"Please reply to this email...\" (block 2).\nawait replaceInParagraph(\n  36,\n  \"Please reply to this email with a copy of these documents so that we have make the necessary arrangements for you, including accommodation and transportation.\",\n  \"Por favor, responda a este e-mail com uma c\\u00F3pia destes documentos para que possamos tomar as provid\\u00EAncias necess\\u00E1rias, incluindo alojamento e transporte.\",\n  1\n);\n\n// Paragraph 37: \"If you have any questions, please contact us via live chat or WhatsApp.\" (block 2).\nawait replaceInParagraph(\n  37,\n  \"If you have any questions, please contact us via \",\n  \"Para mais informa\\u00E7\\u00F5es, contacte-nos atrav\\u00E9s de \",\n  1\n);\nawait replaceInParagraph(37, \" or \", \" ou \", 1);\n\n// Paragraph 38: \"If you have any questions, please contact your country manager, [NAME], at ...\" (block 2).\nawait replaceInParagraph(\n  38,\n  \"If you have any questions, please contact your country manager, \",\n  \"Em caso de d\\u00FAvidas, contacte o seu gestor de parcerias, \",\n  1\n);\nawait replaceInParagraph(38, \", at \", \", em \", 1);\nawait replaceInParagraph(38, \" or \", \" ou \", 1);\n\n// Paragraph 39: closing line (block 2).\nawait replaceInParagraph(\n  39,\n  \"We look forward to seeing you there!\",\n  \"Contamos com a sua presen\\u00E7a!\",\n  1\n);\n\n// Comment text: \"choose either one\" -> \"escolha um deles\".\nconst comments = context.document.getComments();\ncomments.load(\"items\");\nawait context.sync();\ncomments.items[0].content = \"escolha um deles\";\nawait context.sync();\n", "ps1": "# Apply the English -> Portuguese translation edits described by the diff.\n# Strategy: most phrases only occur once with a single correct translation,\n# or occur twice with the SAME translation both times, so a simple\n# \"Find everywhere, Replace All\" on $d.Content is safe for those. A small\n# number of phrases are shared text that must get DIFFERENT translations\n# depending on which half of the document (before/after the page break)\n# they appear in (\", at \" and \"DD Mmm YYYY\"); those are scoped to a Range\n# built from the page-break position so only the correct occurrence changes.\n\n$d = $word.ActiveDocument\n\nfunction Replace-All($rng, [string]$findText, [string]$replaceText) {\n    $f = $rng.Find\n    $f.ClearFormatting()\n    $f.Replacement.ClearFormatting()\n    $f.Text = $findText\n    $f.Replacement.Text = $replaceText\n    $f.Forward = $true\n    $f.Wrap = 0            # wdFindStop - do not wrap around / do not leave the supplied range\n    $f.Format = $false\n    $f.MatchCase = $true\n    $f.MatchWholeWord = $false\n    $f.MatchWildcards = $false\n    $f.MatchSoundsLike = $false\n    $f.MatchAllWordForms = $false\n    $f.Execute($f.Text, $false, $false, $false, $false, $false, $true, 1, $false, $f.Replacement.Text, 2) | Out-Null\n}\n\n# ---- Global-safe replacements (identical translation for every occurrence) ----\n\nReplace-All $d.Content \"English\" \"Ingl\u00eas\"\nReplace-All $d.Content \" / Portuguese / French / Thai / Vietnamese / Spanish\" \" / Portugu\u00eas / Franc\u00eas / Tailand\u00eas / Vietnamita / Espanhol\"\nReplace-All $d.Content \"Brief\" \"Resumo\"\nReplace-All $d.Content \"An email sent to partners in the target country who RSVPed yes but haven\u2019t sent their documents to us. It will be sent via customer.io\" \"E-mail enviado a parceiros no pa\u00eds-alvo que responderam que sim, mas ainda n\u00e3o nos enviaram os documentos. Ser\u00e1 enviado atrav\u00e9s do customer.io\"\nReplace-All $d.Content \"Target audience\" \"P\u00fablico-alvo\"\nReplace-All $d.Content \"Invited partners who haven\u2019t submitted their documents\" \"Parceiros convidados que ainda n\u00e3o enviaram os seus documentos\"\nReplace-All $d.Content \"Subject line\" \"Linha de assunto\"\nReplace-All $d.Content \"[EVENT NAME]\" \"[NOME DO EVENTO]\"\nReplace-All $d.Content \" \u2014 have you submitted your docs?  \" \" - j\u00e1 enviou os seus documentos?  \"\nReplace-All $d.Content \"Don\u2019t forget to send your documents\" \"N\u00e3o se esque\u00e7a de enviar os seus documentos\"\nReplace-All $d.Content \"[PARTNER NAME]\" \"[NOME DO PARCEIRO]\"\nReplace-All $d.Content \"We\u2019re excited to see you at the upcoming \" \"Estamos ansiosos por v\u00ea-lo no pr\u00f3ximo evento \"\nReplace-All $d.Content \"[insert list of documents required]\" \"[inserir lista dos documentos necess\u00e1rios].\"\nReplace-All $d.Content \"We look forward to seeing you there!\" \"Contamos com a sua presen\u00e7a!\"\n\n# Block 1 (before the page break) specific - but these search strings are\n# unique in the whole document anyway, so a document-wide replace is safe.\nReplace-All $d.Content \"Hi \" \"Ol\u00e1 \"\nReplace-All $d.Content \"To confirm your registration, we need the following documents from you by \" \"Para confirmar a sua inscri\u00e7\u00e3o, precisamos que nos envie os seguintes documentos at\u00e9 dia \"\nReplace-All $d.Content \"Please send a copy of these documents to your country manager, \" \"Por favor, envie uma c\u00f3pia destes documentos ao seu gestor de parcerias, \"\nReplace-All $d.Content \" or \" \" ou \"\nReplace-All $d.Content \" (WhatsApp), so that we can make the necessary arrangements for you, including accommodation and transportation.\" \" (WhatsApp), para podermos tomar as medidas necess\u00e1rias, incluindo alojamento e transporte.\"\nReplace-All $d.Content \"If you have any questions, please contact your country manager.\" \"Se tiver alguma d\u00favida, contacte o gestor do seu pa\u00eds.\"\n\n# Block 2 (after the page break) specific - also unique strings, safe to\n# replace document-wide.\nReplace-All $d.Content \"Dear \" \"Ol\u00e1 \"\nReplace-All $d.Content \"To ensure you have the best experience at this event, we need the following documents from you by \" \"De forma a garantir a melhor experi\u00eancia poss\u00edvel neste evento, \u00e9 necess\u00e1rio que nos envie os seguintes documentos at\u00e9 \"\nReplace-All $d.Content \"Please reply to this email with a copy of these documents so that we have make the necessary arrangements for you, including accommodation and transportation.\" \"Por favor, responda a este e-mail com uma c\u00f3pia destes documentos para que possamos tomar as provid\u00eancias necess\u00e1rias, incluindo alojamento e transporte.\"\nReplace-All $d.Content \"If you have any questions, please contact us via \" \"Para mais informa\u00e7\u00f5es, contacte-nos atrav\u00e9s de \"\nReplace-All $d.Content \"If you have any questions, please contact your country manager, \" \"Em caso de d\u00favidas, contacte o seu gestor de parcerias, \"\n\n# ---- Block-scoped replacements (same source phrase, different translation\n#      depending on which half of the document it's in) ----\n# Locate the manual page break that separates block 1 from block 2, then\n# build two Ranges: [0, pageBreakEnd) and [pageBreakEnd, docEnd).\n\n$pbFind = $d.Content.Find\n$pbFind.ClearFormatting()\n$pbFind.Text = \"^m\"\n$pbFind.Forward = $true\n$pbFind.Wrap = 0\n$pbFind.Execute() | Out-Null\n$splitPos = $d.Content.End\n\n$block1 = $d.Range(0, $splitPos)\n$block2 = $d.Range($splitPos, $d.Content.End)\n\nReplace-All $block1 \", at \" \", para \"\nReplace-All $block2 \", at \" \", em \"\n\nReplace-All $block2 \"DD Mmm YYYY\" \"DD Mmm AAAA\"\n\n# ---- Comment text ----\n$comment = $d.Comments.Item(1)\n$comment.Range.Text = \"escolha um deles\"\n"}
